$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "34.631.93"
$ws.Range("E2").Value = "  +1.88%  "

# Row 3
$ws.Range("D3").Value = "1.791.72"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
Set-TextValue "D5" "224.88"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6
Set-TextValue "D6" "0.553"
$ws.Range("E6").Value = "  -0.66%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
Set-TextValue "D8" "32.52"
$ws.Range("E8").Value = "  +5.71%  "

# Row 9
Set-TextValue "D9" "0.283"
$ws.Range("E9").Value = "  +1.62%  "

# Row 10
Set-TextValue "D10" "0.0667"
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("E11").Value = "  +1.40%  "

# Row 12
$ws.Range("D12").Value = "2.052.64"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
Set-TextValue "D13" "11.07"
$ws.Range("E13").Value = "  +10.70%  "

# Row 14
$ws.Range("D14").Value = "1.781.56"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.670.74"
$ws.Range("E15").Value = "  +2.06%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D16" "0.633"
$ws.Range("E16").Value = "  +0.92%  "

# Row 17
Set-TextValue "D17" "4.28"
$ws.Range("E17").Value = "  +1.91%  "

# Row 18
Set-TextValue "D18" "68.98"
$ws.Range("E18").Value = "  +0.69%  "

# Row 19
Set-TextValue "D19" "253.63"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0760"
$ws.Range("E20").Value = "  +2.72%  "

# Row 21
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
Set-TextValue "D22" "10.37"
$ws.Range("E22").Value = "  +0.73%  "

# Row 23
Set-TextValue "D23" "4.22"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
Set-TextValue "D24" "2.13"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
Set-TextValue "D25" "158.53"
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
Set-TextValue "D26" "16.38"
$ws.Range("E26").Value = "  -0.55%  "

# Row 27
Set-TextValue "D27" "7.07"
$ws.Range("E27").Value = "  +1.69%  "

# Row 28
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("E29").Value = "  +0.44%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D30" "0.0517"
$ws.Range("E30").Value = "  +1.01%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "3.76"
$ws.Range("E31").Value = "  -1.05%  "

# Row 32
$ws.Range("E32").Value = "  +0.10%  "

# Row 33
Set-TextValue "D33" "3.58"
$ws.Range("E33").Value = "  +1.32%  "

# Row 34
Set-TextValue "D34" "1.85"
$ws.Range("E34").Value = "  +5.64%  "

# Row 35
$ws.Range("D35").Value = "1.444.78"
$ws.Range("E35").Value = "  -3.45%  "

# Row 36
Set-TextValue "D36" "1.05"
$ws.Range("E36").Value = "  -0.25%  "

# Row 37
Set-TextValue "D37" "0.0188"
$ws.Range("E37").Value = "  +1.32%  "

# Row 38
Set-TextValue "D38" "0.625"
$ws.Range("E38").Value = "  -0.75%  "

# Row 39
Set-TextValue "D39" "82.75"
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("E40").Value = "  +4.71%  "

# Row 41
Set-TextValue "D41" "2.36"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
Set-TextValue "D42" "0.898"
$ws.Range("E42").Value = "  +1.29%  "

# Row 43
$ws.Range("E43").Value = "  -1.02%  "

# Row 44
$ws.Range("E44").Value = "  -0.93%  "

# Row 45
Set-TextValue "D45" "5.94"
$ws.Range("E45").Value = "  +3.42%  "

# Row 46
$ws.Range("E46").Value = "  -1.32%  "

# Row 47
$ws.Range("D47").Value = "1.951.55"
$ws.Range("E47").Value = "  +0.65%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D48" "1.00"
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "103.56"
$ws.Range("E49").Value = "  +5.46%  "

# Row 50
Set-TextValue "D50" "11.87"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51
$ws.Range("E51").Value = "  +6.05%  "
